# "Generate Report for Handback"
#
# This script updates the localization-status workbook to reflect that the
# de-de and zh-cn handback packages have now come back from translation:
#   - status cells flip from "In Translation" to
#     "Handed back: in sync with en-US"
#   - the per-language "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated for both tracked
#     files (66fbcaba...md and 844d606c...md)
#   - the newly populated "Latest Target File" cells become hyperlinks to
#     the same source file the row is about, mirroring column A
#   - a few columns are widened so the new long file names are readable

$wb = $excel.ActiveWorkbook

$urlPrefix = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9213412bba95fc7309c3da2e1fb2330f843e8405/e2e/"
$file1 = "66fbcaba-bdba-465e-a43c-1fb0f1cefd0b.md"
$file2 = "844d606c-9f47-41fd-87a4-dcf166c53008.md"

$hyperlinkColor = 15570276  # BGR packing of RGB(0x64,0x95,0xED) == the workbook's existing HyperLink font color

# ---------------------------------------------------------------------
# 1) Overview sheet: the "In Translation" status cells for both
#    languages/both files are now handed back & in sync, and the
#    zh-cn / de-de columns get wider to fit the longer status text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# 2) zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 -> 66fbcaba...md
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlPrefix + $file1, "", "", $file1)
$wsZhCn.Range("I2").Font.Underline = $true
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor
$wsZhCn.Range("J2").Value = "66fbcaba-bdba-465e-a43c-1fb0f1cefd0b.aeb2a330b995ed2ad2607cb72ed7db0f31b03652.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-22 11:13:48"

# Row 3 -> 844d606c...md
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlPrefix + $file2, "", "", $file2)
$wsZhCn.Range("I3").Font.Underline = $true
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor
$wsZhCn.Range("J3").Value = "844d606c-9f47-41fd-87a4-dcf166c53008.9d6c594c8b094ddf132830ee8fe987129dd4cfa9.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-22 11:13:48"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# 3) de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 -> 66fbcaba...md
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlPrefix + $file1, "", "", $file1)
$wsDeDe.Range("I2").Font.Underline = $true
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor
$wsDeDe.Range("J2").Value = "66fbcaba-bdba-465e-a43c-1fb0f1cefd0b.aeb2a330b995ed2ad2607cb72ed7db0f31b03652.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-22 11:13:55"

# Row 3 -> 844d606c...md
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlPrefix + $file2, "", "", $file2)
$wsDeDe.Range("I3").Font.Underline = $true
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor
$wsDeDe.Range("J3").Value = "844d606c-9f47-41fd-87a4-dcf166c53008.9d6c594c8b094ddf132830ee8fe987129dd4cfa9.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-22 11:13:55"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1666666666667
